$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the header cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 01:52"

# Update Estados Unidos (row 4) totals
$ws.Range("B4").Value = 878974
$ws.Range("C4").Value = 30257
$ws.Range("D4").Value = 85624
$ws.Range("E4").Value = 743596
$ws.Range("F4").Value = 14997
$ws.Range("G4").Value = 2095
$ws.Range("H4").Value = 49754

# Update Noruega (row 41) totals
$ws.Range("B41").Value = 7401
$ws.Range("C41").Value = 63
$ws.Range("E41").Value = 7175
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 194

# Haiti is newly inserted into the sorted country list ahead of Bahamas,
# which pushes Bahamas / Guyana / Islas Caimanes / Sierra Leona each down
# one row (rows 156-160). Haiti and Guyana also receive updated figures.
$ws.Range("A156").Value = "Haiti"
$ws.Range("B156").Value = 72
$ws.Range("C156").Value = 10
$ws.Range("D156").Value = 2
$ws.Range("E156").Value = 65
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 5

$ws.Range("A157").Value = "Bahamas"
$ws.Range("B157").Value = 72
$ws.Range("C157").Value = 7
$ws.Range("D157").Value = 14
$ws.Range("E157").Value = 47
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 2
$ws.Range("H157").Value = 11

$ws.Range("A158").Value = "Guyana"
$ws.Range("B158").Value = 70
$ws.Range("C158").Value = 3
$ws.Range("D158").Value = 9
$ws.Range("E158").Value = 54
$ws.Range("F158").Value = 5
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 7

$ws.Range("A159").Value = "Islas Caimanes"
$ws.Range("B159").Value = 66
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 7
$ws.Range("E159").Value = 58
$ws.Range("F159").Value = 3
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 1

$ws.Range("A160").Value = "Sierra Leona"
$ws.Range("B160").Value = 64
$ws.Range("C160").Value = 3
$ws.Range("D160").Value = 10
$ws.Range("E160").Value = 53
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 1
